# Apply cryptos list update: refresh prices/volumes and shift ranking rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = 'D2'; Value = '27.986.20'; Numeric = $false }
    @{ Cell = 'E2'; Value = '  -0.57%  '; Numeric = $false }
    @{ Cell = 'D3'; Value = '1.857.64'; Numeric = $false }
    @{ Cell = 'E3'; Value = '  -1.27%  '; Numeric = $false }
    @{ Cell = 'D4'; Value = '1.004'; Numeric = $true }
    @{ Cell = 'E4'; Value = '  +0.31%  '; Numeric = $false }
    @{ Cell = 'D5'; Value = '312.15'; Numeric = $true }
    @{ Cell = 'D6'; Value = '1.004'; Numeric = $true }
    @{ Cell = 'E6'; Value = '  +0.29%  '; Numeric = $false }
    @{ Cell = 'D7'; Value = '0.5081'; Numeric = $true }
    @{ Cell = 'E7'; Value = '  +0.28%  '; Numeric = $false }
    @{ Cell = 'D8'; Value = '0.3834'; Numeric = $true }
    @{ Cell = 'E8'; Value = '  -0.62%  '; Numeric = $false }
    @{ Cell = 'D9'; Value = '0.08241'; Numeric = $true }
    @{ Cell = 'E9'; Value = '  -9.05%  '; Numeric = $false }
    @{ Cell = 'E10'; Value = '  -1.52%  '; Numeric = $false }
    @{ Cell = 'B11'; Value = 'Polkadot'; Numeric = $false }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Numeric = $false }
    @{ Cell = 'D11'; Value = '6.201'; Numeric = $true }
    @{ Cell = 'E11'; Value = '  -2.63%  '; Numeric = $false }
    @{ Cell = 'B12'; Value = 'Solana'; Numeric = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Numeric = $false }
    @{ Cell = 'D12'; Value = '20.56'; Numeric = $true }
    @{ Cell = 'E12'; Value = '  -1.24%  '; Numeric = $false }
    @{ Cell = 'B13'; Value = 'WrappedEther'; Numeric = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Numeric = $false }
    @{ Cell = 'D13'; Value = '1.858.26'; Numeric = $false }
    @{ Cell = 'E13'; Value = '  -0.77%  '; Numeric = $false }
    @{ Cell = 'B14'; Value = 'Chainlink'; Numeric = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Numeric = $false }
    @{ Cell = 'D14'; Value = '7.240'; Numeric = $true }
    @{ Cell = 'E14'; Value = '  -0.54%  '; Numeric = $false }
    @{ Cell = 'B15'; Value = 'BinanceUSD'; Numeric = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Numeric = $false }
    @{ Cell = 'D15'; Value = '1.004'; Numeric = $true }
    @{ Cell = 'E15'; Value = '  +0.29%  '; Numeric = $false }
    @{ Cell = 'B16'; Value = 'ShibaInu'; Numeric = $false }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Numeric = $false }
    @{ Cell = 'D16'; Value = '0.00001099'; Numeric = $true }
    @{ Cell = 'E16'; Value = '  -1.49%  '; Numeric = $false }
    @{ Cell = 'B17'; Value = 'Litecoin'; Numeric = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Numeric = $false }
    @{ Cell = 'D17'; Value = '90.70'; Numeric = $true }
    @{ Cell = 'E17'; Value = '  -0.80%  '; Numeric = $false }
    @{ Cell = 'B18'; Value = 'TRON'; Numeric = $false }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Numeric = $false }
    @{ Cell = 'D18'; Value = '0.06644'; Numeric = $true }
    @{ Cell = 'E18'; Value = '  +0.58%  '; Numeric = $false }
    @{ Cell = 'B19'; Value = 'Avalanche'; Numeric = $false }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Numeric = $false }
    @{ Cell = 'D19'; Value = '17.66'; Numeric = $true }
    @{ Cell = 'E19'; Value = '  -3.26%  '; Numeric = $false }
    @{ Cell = 'B20'; Value = 'Dai'; Numeric = $false }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; Numeric = $false }
    @{ Cell = 'D20'; Value = '1.003'; Numeric = $true }
    @{ Cell = 'E20'; Value = '  +0.24%  '; Numeric = $false }
    @{ Cell = 'B21'; Value = 'Uniswap'; Numeric = $false }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Numeric = $false }
    @{ Cell = 'D21'; Value = '6.013'; Numeric = $true }
    @{ Cell = 'E21'; Value = '  -1.94%  '; Numeric = $false }
    @{ Cell = 'B22'; Value = 'WrappedBTC'; Numeric = $false }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Numeric = $false }
    @{ Cell = 'D22'; Value = '28.012.77'; Numeric = $false }
    @{ Cell = 'E22'; Value = '  -0.55%  '; Numeric = $false }
    @{ Cell = 'B23'; Value = 'Cosmos'; Numeric = $false }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Numeric = $false }
    @{ Cell = 'D23'; Value = '11.07'; Numeric = $true }
    @{ Cell = 'E23'; Value = '  -3.58%  '; Numeric = $false }
    @{ Cell = 'B24'; Value = 'Toncoin'; Numeric = $false }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Numeric = $false }
    @{ Cell = 'D24'; Value = '2.243'; Numeric = $true }
    @{ Cell = 'E24'; Value = '  -1.11%  '; Numeric = $false }
    @{ Cell = 'B25'; Value = 'WrappedliquidstakedEther2.0'; Numeric = $false }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Numeric = $false }
    @{ Cell = 'D25'; Value = '2.071.55'; Numeric = $false }
    @{ Cell = 'E25'; Value = '  -0.96%  '; Numeric = $false }
    @{ Cell = 'B26'; Value = 'LidoDAOToken'; Numeric = $false }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Numeric = $false }
    @{ Cell = 'D26'; Value = '2.512'; Numeric = $true }
    @{ Cell = 'E26'; Value = '  -1.36%  '; Numeric = $false }
    @{ Cell = 'B27'; Value = 'Monero'; Numeric = $false }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Numeric = $false }
    @{ Cell = 'D27'; Value = '157.33'; Numeric = $true }
    @{ Cell = 'E27'; Value = '  +0.46%  '; Numeric = $false }
    @{ Cell = 'B28'; Value = 'EthereumClassic'; Numeric = $false }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Numeric = $false }
    @{ Cell = 'D28'; Value = '20.47'; Numeric = $true }
    @{ Cell = 'E28'; Value = '  -1.89%  '; Numeric = $false }
    @{ Cell = 'B29'; Value = 'BitcoinCash'; Numeric = $false }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Numeric = $false }
    @{ Cell = 'D29'; Value = '124.72'; Numeric = $true }
    @{ Cell = 'E29'; Value = '  -1.88%  '; Numeric = $false }
    @{ Cell = 'B30'; Value = 'Stellar'; Numeric = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Numeric = $false }
    @{ Cell = 'D30'; Value = '0.1058'; Numeric = $true }
    @{ Cell = 'E30'; Value = '  -0.30%  '; Numeric = $false }
    @{ Cell = 'B31'; Value = 'ImmutableX'; Numeric = $false }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Numeric = $false }
    @{ Cell = 'D31'; Value = '1.032'; Numeric = $true }
    @{ Cell = 'E31'; Value = '  -3.09%  '; Numeric = $false }
    @{ Cell = 'B32'; Value = 'Filecoin'; Numeric = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Numeric = $false }
    @{ Cell = 'D32'; Value = '5.893'; Numeric = $true }
    @{ Cell = 'E32'; Value = '  +4.76%  '; Numeric = $false }
    @{ Cell = 'B33'; Value = 'HuobiToken'; Numeric = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Numeric = $false }
    @{ Cell = 'D33'; Value = '3.600'; Numeric = $true }
    @{ Cell = 'E33'; Value = '  +0.05%  '; Numeric = $false }
    @{ Cell = 'B34'; Value = 'FraxShare'; Numeric = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Numeric = $false }
    @{ Cell = 'D34'; Value = '9.370'; Numeric = $true }
    @{ Cell = 'E34'; Value = '  -2.46%  '; Numeric = $false }
    @{ Cell = 'B35'; Value = 'VeChain'; Numeric = $false }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Numeric = $false }
    @{ Cell = 'D35'; Value = '0.02410'; Numeric = $true }
    @{ Cell = 'E35'; Value = '  +0.01%  '; Numeric = $false }
    @{ Cell = 'D36'; Value = '0.06500'; Numeric = $true }
    @{ Cell = 'E36'; Value = '  -1.99%  '; Numeric = $false }
    @{ Cell = 'B37'; Value = 'Algorand'; Numeric = $false }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Numeric = $false }
    @{ Cell = 'D37'; Value = '0.2172'; Numeric = $true }
    @{ Cell = 'E37'; Value = '  -1.06%  '; Numeric = $false }
    @{ Cell = 'B38'; Value = 'TheSandbox'; Numeric = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Numeric = $false }
    @{ Cell = 'D38'; Value = '0.6535'; Numeric = $true }
    @{ Cell = 'E38'; Value = '  +1.51%  '; Numeric = $false }
    @{ Cell = 'B39'; Value = 'ARBITRUM'; Numeric = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Numeric = $false }
    @{ Cell = 'D39'; Value = '1.197'; Numeric = $true }
    @{ Cell = 'E39'; Value = '  -1.54%  '; Numeric = $false }
    @{ Cell = 'B40'; Value = 'InternetComputer(DFINITY)'; Numeric = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Numeric = $false }
    @{ Cell = 'D40'; Value = '4.993'; Numeric = $true }
    @{ Cell = 'E40'; Value = '  +1.20%  '; Numeric = $false }
    @{ Cell = 'B41'; Value = 'TrustWalletToken'; Numeric = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Numeric = $false }
    @{ Cell = 'D41'; Value = '1.219'; Numeric = $true }
    @{ Cell = 'E41'; Value = '  -5.33%  '; Numeric = $false }
    @{ Cell = 'B42'; Value = 'Aptos'; Numeric = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Numeric = $false }
    @{ Cell = 'D42'; Value = '11.17'; Numeric = $true }
    @{ Cell = 'E42'; Value = '  -3.05%  '; Numeric = $false }
    @{ Cell = 'B43'; Value = 'Decentraland'; Numeric = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; Numeric = $false }
    @{ Cell = 'D43'; Value = '0.6123'; Numeric = $true }
    @{ Cell = 'E43'; Value = '  +1.12%  '; Numeric = $false }
    @{ Cell = 'B44'; Value = 'EnergySwap'; Numeric = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Numeric = $false }
    @{ Cell = 'D44'; Value = '13.07'; Numeric = $true }
    @{ Cell = 'E44'; Value = '  -1.80%  '; Numeric = $false }
    @{ Cell = 'B45'; Value = 'WEMIXTOKEN'; Numeric = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; Numeric = $false }
    @{ Cell = 'D45'; Value = '1.282'; Numeric = $true }
    @{ Cell = 'E45'; Value = '  +0.49%  '; Numeric = $false }
    @{ Cell = 'B46'; Value = 'PancakeSwap'; Numeric = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Numeric = $false }
    @{ Cell = 'D46'; Value = '3.650'; Numeric = $true }
    @{ Cell = 'E46'; Value = '  -0.52%  '; Numeric = $false }
    @{ Cell = 'B47'; Value = 'NEARProtocol'; Numeric = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Numeric = $false }
    @{ Cell = 'D47'; Value = '2.011'; Numeric = $true }
    @{ Cell = 'E47'; Value = '  +0.19%  '; Numeric = $false }
    @{ Cell = 'B48'; Value = 'EOS'; Numeric = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; Numeric = $false }
    @{ Cell = 'D48'; Value = '1.205'; Numeric = $true }
    @{ Cell = 'E48'; Value = '  -3.07%  '; Numeric = $false }
    @{ Cell = 'B49'; Value = 'Quant'; Numeric = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Numeric = $false }
    @{ Cell = 'D49'; Value = '119.95'; Numeric = $true }
    @{ Cell = 'E49'; Value = '  -1.18%  '; Numeric = $false }
    @{ Cell = 'B50'; Value = 'Aave'; Numeric = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; Numeric = $false }
    @{ Cell = 'D50'; Value = '78.41'; Numeric = $true }
    @{ Cell = 'E50'; Value = '  -1.51%  '; Numeric = $false }
    @{ Cell = 'B51'; Value = 'Cronos'; Numeric = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Numeric = $false }
    @{ Cell = 'D51'; Value = '0.06832'; Numeric = $true }
    @{ Cell = 'E51'; Value = '  -1.47%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $r.NumberFormat = '@'
        $r.Value = $u.Value
        $r.Style = 'Normal'
    } else {
        $r.Value = $u.Value
    }
}
